$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Text-like values that must NOT be auto-converted (dates, etc.) get a
# leading apostrophe so Excel stores them as literal text; the trailing
# Style reset below strips the quote-prefix formatting back to Normal so
# no stray number-format style is introduced.
$ws.Cells.Item($row, 1).Value  = "DF"
$ws.Cells.Item($row, 2).Value  = "TES1656"
$ws.Cells.Item($row, 3).Value  = ""
$ws.Cells.Item($row, 4).Value  = ""
$ws.Cells.Item($row, 5).Value  = ""
$ws.Cells.Item($row, 6).Value  = "T"
$ws.Cells.Item($row, 7).Value  = "T"
$ws.Cells.Item($row, 8).Value  = "T - (T 03/11/25_24H) - DF"
$ws.Cells.Item($row, 9).Value  = "'03/11/25"
$ws.Cells.Item($row, 10).Value = "24H"
$ws.Cells.Item($row, 11).Value = "14/11/25"
$ws.Cells.Item($row, 12).Value = "DENTRO"
$ws.Cells.Item($row, 13).Value = ""

# Ensure every cell in the new row (including the blank ones, and the
# quote-prefixed date text cell) gets written out / normalized back to
# the default "Normal" style, matching the rest of the sheet.
for ($col = 1; $col -le 13; $col++) {
    $ws.Cells.Item($row, $col).Style = "Normal"
}
